$p = $ppt.ActivePresentation

# Locate the "VIDEO LINKS" slide (slide 12 in this deck).
$s = $p.Slides.Item(12)

# --- Title shape: "VIDEO LINKS" -> "LINKS" ---
$titleShape = $s.Shapes.Title
$titleShape.TextFrame.TextRange.Text = "LINKS"

# --- TextBox 6: rebuild the links list with hyperlinks ---
$linksShape = $s.Shapes.Item("TextBox 6")
$tr = $linksShape.TextFrame.TextRange

$dash = [string][char]0x2013

$repoLabel = "Team Repository " + $dash + " "
$repoUrl   = "https://github.com/bensonnd/DS6306_David_Neil_Project_1"

$videoLabel = "Neil YouTube Video - "
$videoUrl   = "https://youtu.be/s5X5uf_ModA "

# Paragraph 1: repo label + link, Paragraph 2: blank, Paragraph 3: video label + link
$tr.Text = $repoLabel + $repoUrl + "`r`r" + $videoLabel + $videoUrl

# Hyperlink the repo URL text in paragraph 1
$repoPara = $tr.Paragraphs(1, 1)
$repoUrlRange = $repoPara.Characters($repoLabel.Length + 1, $repoUrl.Length)
$repoUrlRange.ActionSettings.Item(1).Hyperlink.Address = $repoUrl

# Hyperlink the video URL text in paragraph 3
$videoPara = $tr.Paragraphs(3, 1)
$videoUrlRange = $videoPara.Characters($videoLabel.Length + 1, $videoUrl.Length)
$videoUrlRange.ActionSettings.Item(1).Hyperlink.Address = $videoUrl.Trim()
